$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 1930.4348
$ws.Range("I113").Value = 1650
$ws.Range("J113").Value = 1989.4736
$ws.Range("K113").Value = 1650
$ws.Range("L113").Value = 1989.4736
$ws.Range("M113").Value = 1604
$ws.Range("N113").Value = -8497.473599999999
$ws.Range("H116").Value = 6414924.5
$ws.Range("I116").Value = 12827249
$ws.Range("J116").Value = 2600
$ws.Range("K116").Value = 12827249
$ws.Range("L116").Value = 2600
$ws.Range("M116").Value = -12823807
$ws.Range("N116").Value = -9484
$ws.Range("H132").Value = 4611.552
$ws.Range("I132").Value = 4503.7617
$ws.Range("J132").Value = 4894.5
$ws.Range("K132").Value = 13511.2851
$ws.Range("L132").Value = 14683.5
$ws.Range("M132").Value = -10981.2851
$ws.Range("N132").Value = -19743.5
$ws.Range("H137").Value = 38135.855
$ws.Range("I137").Value = 1665
$ws.Range("J137").Value = 94499.91
$ws.Range("K137").Value = 4995
$ws.Range("L137").Value = 283499.73
$ws.Range("M137").Value = -2445
$ws.Range("N137").Value = -288599.73

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1687.3846
$ws.Range("I61").Value = 1248.1852
$ws.Range("J61").Value = 2675.5833
$ws.Range("K61").Value = 1248.1852
$ws.Range("L61").Value = 2675.5833
$ws.Range("M61").Value = -1036.1852
$ws.Range("N61").Value = -3099.5833
$ws.Range("H74").Value = 2985.2654
$ws.Range("I74").Value = 3688.6667
$ws.Range("J74").Value = 1037.3846
$ws.Range("K74").Value = 3688.6667
$ws.Range("L74").Value = 1037.3846
$ws.Range("M74").Value = -2814.6667
$ws.Range("N74").Value = -2785.3846
$ws.Range("H77").Value = 2985.2654
$ws.Range("I77").Value = 3688.6667
$ws.Range("J77").Value = 1037.3846
$ws.Range("K77").Value = 18443.3335
$ws.Range("L77").Value = 5186.923000000001
$ws.Range("M77").Value = -14075.3335
$ws.Range("N77").Value = -13922.923
$ws.Range("H122").Value = 1185.8889
$ws.Range("I122").Value = 1146.625
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 3439.875
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -989.875
$ws.Range("N122").Value = -9400
$ws.Range("H132").Value = 2119.577
$ws.Range("I132").Value = 1956.1177
$ws.Range("J132").Value = 2428.3333
$ws.Range("K132").Value = 5868.3531
$ws.Range("L132").Value = 7284.999899999999
$ws.Range("M132").Value = -3338.3531
$ws.Range("N132").Value = -12344.9999
$ws.Range("H136").Value = 1687.3846
$ws.Range("I136").Value = 1248.1852
$ws.Range("J136").Value = 2675.5833
$ws.Range("K136").Value = 3744.5556
$ws.Range("L136").Value = 8026.749899999999
$ws.Range("M136").Value = -1194.5556
$ws.Range("N136").Value = -13126.7499
$ws.Range("H139").Value = 29091.25
$ws.Range("J139").Value = 32571.666
$ws.Range("L139").Value = 32571.666
$ws.Range("N139").Value = -42851.666

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 494.27274
$ws.Range("I94").Value = 446.33334
$ws.Range("J94").Value = 551.8
$ws.Range("K94").Value = 446.33334
$ws.Range("L94").Value = 551.8
$ws.Range("M94").Value = 4.666659999999979
$ws.Range("N94").Value = -1453.8
$ws.Range("H99").Value = 2102.889
$ws.Range("I99").Value = 1063
$ws.Range("K99").Value = 1063
$ws.Range("M99").Value = 435
$ws.Range("H107").Value = 1961.6207
$ws.Range("I107").Value = 2045.0416
$ws.Range("K107").Value = 2045.0416
$ws.Range("M107").Value = -125.0416
$ws.Range("H134").Value = 1697.7925
$ws.Range("I134").Value = 1471.0465
$ws.Range("J134").Value = 2672.8
$ws.Range("K134").Value = 4413.139499999999
$ws.Range("L134").Value = 8018.400000000001
$ws.Range("M134").Value = -1878.139499999999
$ws.Range("N134").Value = -13088.4

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3461.0908
$ws.Range("I31").Value = 2042.421
$ws.Range("K31").Value = 2042.421
$ws.Range("M31").Value = -1747.421
$ws.Range("H34").Value = 3461.0908
$ws.Range("I34").Value = 2042.421
$ws.Range("K34").Value = 2042.421
$ws.Range("M34").Value = -1840.421
$ws.Range("H99").Value = 1933.8462
$ws.Range("I99").Value = 1605.3529
$ws.Range("J99").Value = 2187.682
$ws.Range("K99").Value = 1605.3529
$ws.Range("L99").Value = 2187.682
$ws.Range("M99").Value = -107.3529000000001
$ws.Range("N99").Value = -5183.682
$ws.Range("H126").Value = 1933.8462
$ws.Range("I126").Value = 1605.3529
$ws.Range("J126").Value = 2187.682
$ws.Range("K126").Value = 4816.0587
$ws.Range("L126").Value = 6563.045999999999
$ws.Range("M126").Value = -2346.0587
$ws.Range("N126").Value = -11503.046
$ws.Range("H132").Value = 1967.44
$ws.Range("I132").Value = 2045.2307
$ws.Range("J132").Value = 1883.1666
$ws.Range("K132").Value = 6135.6921
$ws.Range("L132").Value = 5649.4998
$ws.Range("M132").Value = -3605.6921
$ws.Range("N132").Value = -10709.4998
$ws.Range("H134").Value = 2276.2896
$ws.Range("I134").Value = 1370.069
$ws.Range("J134").Value = 5196.3335
$ws.Range("K134").Value = 4110.207
$ws.Range("L134").Value = 15589.0005
$ws.Range("M134").Value = -1575.207
$ws.Range("N134").Value = -20659.0005

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 615.8077
$ws.Range("I5").Value = 565.7826
$ws.Range("K5").Value = 1697.3478
$ws.Range("M5").Value = -1585.3478
$ws.Range("H135").Value = 615.8077
$ws.Range("I135").Value = 565.7826
$ws.Range("K135").Value = 5092.0434
$ws.Range("M135").Value = -2557.0434
$ws.Range("H137").Value = 3139.2068
$ws.Range("I137").Value = 2973.7693
$ws.Range("J137").Value = 3273.625
$ws.Range("K137").Value = 8921.3079
$ws.Range("L137").Value = 9820.875
$ws.Range("M137").Value = -3821.3079
$ws.Range("N137").Value = -20020.875

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 7144388.5
$ws.Range("I122").Value = 9092331
$ws.Range("J122").Value = 1933.3334
$ws.Range("K122").Value = 27276993
$ws.Range("L122").Value = 5800.0002
$ws.Range("M122").Value = -27274543
$ws.Range("N122").Value = -10700.0002
$ws.Range("H132").Value = 2845.7925
$ws.Range("I132").Value = 2744.3901
$ws.Range("J132").Value = 3192.25
$ws.Range("K132").Value = 8233.1703
$ws.Range("L132").Value = 9576.75
$ws.Range("M132").Value = -5703.1703
$ws.Range("N132").Value = -14636.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3000.5
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").Value = ""
$ws.Range("H71").Value = 3000.5
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = ""
$ws.Range("H93").Value = 18551
$ws.Range("I93").Value = 26600.5
$ws.Range("J93").Value = 2452
$ws.Range("K93").Value = 26600.5
$ws.Range("L93").Value = 2452
$ws.Range("M93").Value = -25352.5
$ws.Range("N93").Value = -4948
$ws.Range("H122").Value = 4773.25
$ws.Range("J122").Value = 5221.8
$ws.Range("L122").Value = 15665.4
$ws.Range("N122").Value = -20565.4
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = ""

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").Value = ""
$ws.Range("H122").Value = 2415.9678
$ws.Range("I122").Value = 2224.2632
$ws.Range("J122").Value = 2719.5
$ws.Range("K122").Value = 6672.7896
$ws.Range("L122").Value = 8158.5
$ws.Range("M122").Value = -4222.7896
$ws.Range("N122").Value = -13058.5
$ws.Range("H132").Value = 2731.907
$ws.Range("I132").Value = 2987.5417
$ws.Range("J132").Value = 2409
$ws.Range("K132").Value = 8962.625100000001
$ws.Range("L132").Value = 7227
$ws.Range("M132").Value = -6432.625100000001
$ws.Range("N132").Value = -12287
